$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 93, shifting rows 93:143 down to 94:144
$ws.Rows.Item(93).Insert()

# Fill in the new row 93 with data (copy constant columns from what is now row 94,
# i.e. the old row 93 content, then set the changed values)
$ws.Range("A93").Value = 5
$ws.Range("B93").Value = "Macroferia Regional de Talca"
$ws.Range("C93").Value = "Maule"
$ws.Range("D93").Value = 44455
$ws.Range("E93").Value = 7
$ws.Range("F93").Value = "Fruta"
$ws.Range("G93").Value = 100108
$ws.Range("H93").Value = "Tropicales y subtropicales"
$ws.Range("I93").Value = 100108005
$ws.Range("J93").Value = "Piña"
$ws.Range("K93").Value = "Caramelo"
$ws.Range("L93").Value = "Segunda"
$ws.Range("M93").Value = 54
$ws.Range("N93").Value = 22000
$ws.Range("O93").Value = 22000
$ws.Range("P93").Value = 22000
$ws.Range("Q93").Value = "$/caja 14 unidades"
$ws.Range("R93").Value = "Ecuador"
$ws.Range("S93").Value = 1571
$ws.Range("T93").Value = 14
